# Collapse the four "2018年キャンペーン期間 (対象：ペルセウス)：、10月30日〜11月8日、
# 11月29日〜12月8日" runs-per-character blocks (one per language/size copy of the
# flyer) down to a single plain run holding the already-translated date string,
# exactly as the commit describes: "Insert traducted dates ... constellation's
# name must be traducted in the next step".
#
# The paragraph keeps its <w:pPr> (and any <w:sectPr> inside it) untouched; only
# the run content changes, and the replacement run carries no <w:rPr> (it simply
# inherits the paragraph/style defaults), matching the target XML:
#   <w:r><w:t>年キャンペーン期間 対象：Pegasus: 10月8〜17日、11月7〜16日、</w:t></w:r>

$d = $word.ActiveDocument

$oldText = "2018年キャンペーン期間 (対象：ペルセウス)：、10月30日〜11月8日、11月29日〜12月8日"
$newText = "年キャンペーン期間 対象：Pegasus: 10月8〜17日、11月7〜16日、"

$packageXml = "<?xml version=`"1.0`" standalone=`"yes`"?>" + `
  "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" + `
    "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" + `
      "<pkg:xmlData>" + `
        "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" + `
          "<w:body><w:p><w:r><w:t>" + $newText + "</w:t></w:r></w:p></w:body>" + `
        "</w:document>" + `
      "</pkg:xmlData>" + `
    "</pkg:part>" + `
  "</pkg:package>"

$searchRange = $d.Content
$replaced = 0

$found = $searchRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found -and $replaced -lt 20) {
    $replaced = $replaced + 1

    # Use a freshly-minted Range over the match bounds for the edit itself -
    # InsertXML replaces that range's contents with the runs found in the
    # fragment's <w:body>, dropping all of the original runs/proofErr marks
    # and leaving a single unformatted run behind.
    $editRange = $d.Range($searchRange.Start, $searchRange.End)
    $editRange.InsertXML($packageXml)

    # Keep searching the remainder of the document for further occurrences.
    $searchRange = $d.Range($editRange.End, $d.Content.End)
    $found = $searchRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

Write-Output "Replaced $replaced occurrence(s)"
